$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 15-45 (the old, expanded per-field rows) so the sheet
# collapses down to the new condensed A1:A14 layout.
$ws.Range("A15:A45").EntireRow.Delete()

# Row 1 (the set title) is unchanged. Rows 2-14 now hold one
# Python-tuple-repr string per token instead of one field per row.
$ws.Range("A2").Value = @'
('Beast', ['Token Creature — Beast', 'Trample', '4/4'])
'@

$ws.Range("A3").Value = @'
('Centaur', ['Token Creature — Centaur', '3/3'])
'@

$ws.Range("A4").Value = @'
('Domri, Chaos Bringer Emblem', ['Emblem — Domri', 'At the beginning of each end step, create a 4/4 red and green Beast creature token with trample.'])
'@

$ws.Range("A5").Value = @'
('Frog Lizard', ['Token Creature — Frog Lizard', '3/3'])
'@

$ws.Range("A6").Value = @'
('Goblin', ['Token Creature — Goblin', '1/1'])
'@

$ws.Range("A7").Value = @'
('Human', ['Token Creature — Human', '1/1'])
'@

$ws.Range("A8").Value = @'
('Illusion', ['Token Creature — Illusion', 'Whenever this creature blocks a creature, that creature doesn’t untap during its controller’s next untap step.', '0/2'])
'@

$ws.Range("A9").Value = @'
('Ooze', ['Token Creature — Ooze', '2/2'])
'@

$ws.Range("A10").Value = @'
('Sphinx', ['Token Creature — Sphinx', 'Flying, vigilance', '4/4'])
'@

$ws.Range("A11").Value = @'
('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])
'@

$ws.Range("A12").Value = @'
('Thopter', ['Token Artifact Creature — Thopter', 'Flying', '1/1'])
'@

$ws.Range("A13").Value = @'
('Treasure', ['Token Artifact — Treasure', '{T}, Sacrifice this artifact: Add one mana of any color.'])
'@

$ws.Range("A14").Value = @'
('Zombie', ['Token Creature — Zombie', '2/2'])
'@
